# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.471.64"
$ws.Range("E2").Value = "  +6.49%  "

# Row 3
$ws.Range("D3").Value = "'1.718.39"
$ws.Range("E3").Value = "  +3.16%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'333.26"
$ws.Range("E5").Value = "  +1.20%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").Value = "'0.3703"
$ws.Range("E7").Value = "  +1.52%  "

# Row 8
$ws.Range("D8").Value = "'48.18"
$ws.Range("E8").Value = "  +1.62%  "

# Row 9
$ws.Range("D9").Value = "'0.3353"
$ws.Range("E9").Value = "  +2.44%  "

# Row 10
$ws.Range("D10").Value = "'1.183"
$ws.Range("E10").Value = "  +3.61%  "

# Row 11
$ws.Range("D11").Value = "'0.07372"
$ws.Range("E11").Value = "  +3.89%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.04%  "

# Row 13
$ws.Range("D13").Value = "'6.361"
$ws.Range("E13").Value = "  +4.32%  "

# Row 14
$ws.Range("E14").Value = "  +2.38%  "

# Row 15
$ws.Range("D15").Value = "'7.024"
$ws.Range("E15").Value = "  +6.02%  "

# Row 16
$ws.Range("D16").Value = "'1.717.48"
$ws.Range("E16").Value = "  +3.38%  "

# Row 17
$ws.Range("E17").Value = "  +1.55%  "

# Row 18
$ws.Range("D18").Value = "'0.06625"
$ws.Range("E18").Value = "  -0.63%  "

# Row 19
$ws.Range("D19").Value = "'81.80"
$ws.Range("E19").Value = "  +3.89%  "

# Row 20
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("D21").Value = "'16.49"
$ws.Range("E21").Value = "  +4.02%  "

# Row 22
$ws.Range("D22").Value = "'6.114"
$ws.Range("E22").Value = "  +2.89%  "

# Row 23
$ws.Range("D23").Value = "'12.73"
$ws.Range("E23").Value = "  +1.08%  "

# Row 24
$ws.Range("D24").Value = "'26.453.52"
$ws.Range("E24").Value = "  +6.47%  "

# Row 25
$ws.Range("D25").Value = "'2.436"
$ws.Range("E25").Value = "  -1.48%  "

# Row 26
$ws.Range("D26").Value = "'2.379"
$ws.Range("E26").Value = "  -2.72%  "

# Row 27
$ws.Range("D27").Value = "'1.382"
$ws.Range("E27").Value = "  +17.93%  "

# Row 28
$ws.Range("D28").Value = "'151.32"
$ws.Range("E28").Value = "  +0.97%  "

# Row 29
$ws.Range("D29").Value = "'19.32"
$ws.Range("E29").Value = "  +3.08%  "

# Row 30
$ws.Range("D30").Value = "'1.911.78"
$ws.Range("E30").Value = "  +3.67%  "

# Row 31
$ws.Range("D31").Value = "'130.86"
$ws.Range("E31").Value = "  +3.67%  "

# Row 32
$ws.Range("E32").Value = "  +1.02%  "

# Row 33
$ws.Range("D33").Value = "'5.890"
$ws.Range("E33").Value = "  +3.01%  "

# Row 34
$ws.Range("D34").Value = "'0.08607"
$ws.Range("E34").Value = "  +1.43%  "

# Row 35
$ws.Range("D35").Value = "'1.703"
$ws.Range("E35").Value = "  +3.16%  "

# Row 36
$ws.Range("D36").Value = "'12.62"
$ws.Range("E36").Value = "  +3.28%  "

# Row 37
$ws.Range("D37").Value = "'5.335"
$ws.Range("E37").Value = "  +2.92%  "

# Row 38
$ws.Range("D38").Value = "'0.02310"
$ws.Range("E38").Value = "  +1.53%  "

# Row 39
$ws.Range("E39").Value = "  +3.07%  "

# Row 40
$ws.Range("D40").Value = "'0.06180"
$ws.Range("E40").Value = "  -0.65%  "

# Row 41
$ws.Range("D41").Value = "'8.380"
$ws.Range("E41").Value = "  +1.63%  "

# Row 42
$ws.Range("D42").Value = "'1.218"
$ws.Range("E42").Value = "  -1.91%  "

# Row 43
$ws.Range("D43").Value = "'0.6155"
$ws.Range("E43").Value = "  +3.51%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.15"
$ws.Range("E44").Value = "  +4.66%  "

# Row 45
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("D46").Value = "'3.902"
$ws.Range("E46").Value = "  +1.35%  "

# Row 47
$ws.Range("D47").Value = "'0.5953"
$ws.Range("E47").Value = "  +5.07%  "

# Row 48
$ws.Range("D48").Value = "'128.18"
$ws.Range("E48").Value = "  +2.02%  "

# Row 49
$ws.Range("D49").Value = "'2.032"
$ws.Range("E49").Value = "  +3.60%  "

# Row 50
$ws.Range("D50").Value = "'0.07153"
$ws.Range("E50").Value = "  +2.28%  "

# Row 51
$ws.Range("D51").Value = "'76.66"
$ws.Range("E51").Value = "  +1.62%  "

Write-Host "Updated cryptos list"